$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.525.77'
$ws.Range("E2").Value = '  +0.13%  '

$ws.Range("D3").Value = '1.794.88'
$ws.Range("E3").Value = '  -0.92%  '

$ws.Range("D4").NumberFormat = '@'
$ws.Range("D4").Value = '0.9996'
$ws.Range("E4").Value = '  -0.83%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '339.58'
$ws.Range("E5").Value = '  +1.78%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '0.9982'
$ws.Range("E6").Value = '  -0.46%  '

$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.3920'
$ws.Range("E7").Value = '  +2.96%  '

$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.3466'
$ws.Range("E8").Value = '  -1.10%  '

$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '48.34'
$ws.Range("E9").Value = '  -1.67%  '

$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '1.200'
$ws.Range("E10").Value = '  -1.72%  '

$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.07507'
$ws.Range("E11").Value = '  -1.98%  '

$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '0.9971'
$ws.Range("E12").Value = '  -0.88%  '

$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '21.93'
$ws.Range("E13").Value = '  -0.83%  '

$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '6.520'
$ws.Range("E14").Value = '  -0.47%  '

$ws.Range("D15").Value = '1.792.87'
$ws.Range("E15").Value = '  -1.46%  '

$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '7.175'
$ws.Range("E16").Value = '  +0.80%  '

$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '0.00001102'
$ws.Range("E17").Value = '  -0.81%  '

$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '0.06693'
$ws.Range("E18").Value = '  -0.02%  '

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '84.99'
$ws.Range("E19").Value = '  -1.13%  '

$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '0.9968'
$ws.Range("E20").Value = '  -0.60%  '

$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '17.74'
$ws.Range("E21").Value = '  +1.62%  '

$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '6.573'
$ws.Range("E22").Value = '  +1.25%  '

$ws.Range("D23").Value = '27.507.57'
$ws.Range("E23").Value = '  -0.06%  '

$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '12.45'
$ws.Range("E24").Value = '  -4.05%  '

$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '2.411'
$ws.Range("E25").Value = '  -1.34%  '

$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '2.514'
$ws.Range("E26").Value = '  -4.04%  '

$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '21.23'
$ws.Range("E27").Value = '  -2.86%  '

$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '1.465'
$ws.Range("E28").Value = '  +0.54%  '

$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '156.66'
$ws.Range("E29").Value = '  +3.97%  '

$ws.Range("D30").Value = '1.999.43'
$ws.Range("E30").Value = '  -1.09%  '

$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '135.35'
$ws.Range("E31").Value = '  +0.66%  '

$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '4.046'
$ws.Range("E32").Value = '  -1.00%  '

$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '6.058'
$ws.Range("E33").Value = '  -1.77%  '

$ws.Range("E34").Value = '  +0.37%  '

$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '13.08'
$ws.Range("E35").Value = '  -3.53%  '

$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '1.625'
$ws.Range("E36").Value = '  -3.40%  '

$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '5.459'
$ws.Range("E37").Value = '  -1.39%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '0.02413'
$ws.Range("E38").Value = '  +1.85%  '

$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '0.06509'
$ws.Range("E39").Value = '  +1.26%  '

$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '0.6826'
$ws.Range("E40").Value = '  -1.30%  '

$ws.Range("E41").Value = '  -0.85%  '

$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '1.252'
$ws.Range("E42").Value = '  -3.12%  '

$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '8.402'
$ws.Range("E43").Value = '  -6.70%  '

$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '14.46'
$ws.Range("E44").Value = '  -0.89%  '

$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '0.9971'
$ws.Range("E45").Value = '  -0.49%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '0.6403'
$ws.Range("E46").Value = '  -1.32%  '

$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '3.875'
$ws.Range("E47").Value = '  +0.30%  '

$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '2.146'
$ws.Range("E48").Value = '  -0.32%  '

$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '132.10'
$ws.Range("E49").Value = '  +0.23%  '

$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '0.07191'
$ws.Range("E50").Value = '  -1.24%  '

$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '79.80'
$ws.Range("E51").Value = '  -0.62%  '
